$d = $word.ActiveDocument

# 1) Replace the text of the "14th July 2022" paragraph with the new date.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Testing on 14th July 2022*") {
        $p.Range.Text = "Testing on 20th of February 2023"
        break
    }
}

# 2) Remove the whole "Testing on 19th July 2022" paragraph (text + its paragraph mark).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Testing on 19th July 2022*") {
        $p.Range.Delete()
        break
    }
}
